# Append new sensor-log rows to the SeniorConnect master log workbook.
# Each sheet logs Date/Timestamp/Hour/Location/Value/Status columns as plain
# text (the source system always writes inlineStr/text cells - dates,
# times and percentages must NOT be auto-coerced into Excel date/number
# serials by COM's smart-typing). A leading "'" forces text for the values
# Excel would otherwise re-interpret (the bare date "2026-01-28" and the
# "NN.N%" humidity readings); plain "HH:MM[:SS]" strings and "NN.NC"
# temperature readings already round-trip as text so are left unprefixed.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 79-89 (Bathroom / No Motion / Inactive) ---
$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    @("2026-01-28","16:43:32","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:43:37","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:43:42","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:43:47","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:43:52","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:43:57","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:44:02","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:44:08","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:44:12","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:44:18","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:44:23","16:00","Bathroom","No Motion","Inactive")
)
$r = 79
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Humidity sheet: rows 77-87 (Bathroom / NN.N% / Active) ---
$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    @("2026-01-28","16:43:32","16:00","Bathroom","86.9%","Active"),
    @("2026-01-28","16:43:36","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:43:40","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:43:44","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:43:48","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:43:52","16:00","Bathroom","86.9%","Active"),
    @("2026-01-28","16:43:56","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:44:04","16:00","Bathroom","86.9%","Active"),
    @("2026-01-28","16:44:08","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:44:16","16:00","Bathroom","87.7%","Active"),
    @("2026-01-28","16:44:24","16:00","Bathroom","87.8%","Active")
)
$r = 77
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Temperature sheet: rows 77-87 (Bathroom / NN.NC / Active) ---
$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
    @("2026-01-28","16:43:32","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:43:36","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:43:40","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:43:44","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:43:48","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:43:52","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:43:56","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:44:04","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:44:08","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:44:16","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:44:24","16:00","Bathroom","22.9C","Active")
)
$r = 77
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Proximity sheet: rows 4-5 (Bathroom Door ENTER/EXIT) ---
$ws = $wb.Worksheets.Item("Proximity")
$rows = @(
    @("2026-01-28","16:44:17","16:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-28","16:44:26","16:00","Bathroom Door","EXIT","User EXITED Bathroom")
)
$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- mmWave sheet: rows 3-5 (Living Room / NO_PRESENCE) ---
$ws = $wb.Worksheets.Item("mmWave")
$rows = @(
    @("2026-01-28","16:44:07","16:00","Living Room","NO_PRESENCE","Inactive"),
    @("2026-01-28","16:44:11","16:00","Living Room","NO_PRESENCE","Active"),
    @("2026-01-28","16:44:17","16:00","Living Room","NO_PRESENCE","Inactive")
)
$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
